# Insert a new data record above row 102 (shifts existing rows 102-176 down to 103-177)
# then populate the new row 102 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(102).Insert()

$ws.Range("A102").Value = 4
$ws.Range("B102").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C102").Value = "Los Lagos"
$ws.Range("D102").Value = 44574
$ws.Range("E102").Value = 10
$ws.Range("F102").Value = 100112032
$ws.Range("G102").Value = "Zapallo italiano"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 140
$ws.Range("K102").Value = 15000
$ws.Range("L102").Value = 15000
$ws.Range("M102").Value = 15000
$ws.Range("N102").Value = "$/caja 60 unidades"
$ws.Range("O102").Value = "Región Metropolitana"
$ws.Range("P102").Value = 250
$ws.Range("Q102").Value = 60
$ws.Range("R102").Value = "Hortaliza"
